$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column widths to match the added columns D:K (15 chars wide, same as B/C)
$ws.Range("D1:K1").ColumnWidth = 14.1

# Row 1: File Name
$ws.Range("B1").Value = "Data collected @ Sat Jul 15 10:43:16 2023.txt"
$ws.Range("C1").Value = "Data collected @ Sat Jul 15 11:06:46 2023.txt"
$ws.Range("D1").Value = "Data collected @ Sat Aug 3 8:10:15 2023.txt"
$ws.Range("E1").Value = "Data collected @ Sat Jul 89 10:55:15 2023.txt"
$ws.Range("F1").Value = "Data collected @ Sat Jul 15 10:45:05 2023.txt"
$ws.Range("G1").Value = "Data collected @ Sat Jul 15 10:55:15 2023.txt"

# Row 2: Time Initial
$ws.Range("B2").Value = "10:43:21"
$ws.Range("C2").Value = "11:07:01"
$ws.Range("D2").Value = "9:55:28"
$ws.Range("E2").Value = "8:44:22"
$ws.Range("F2").Value = "10:45:10"
$ws.Range("G2").Value = "10:55:22"

# Row 3: Time Final
$ws.Range("B3").Value = "10:43:51"
$ws.Range("C3").Value = "11:11:39"
$ws.Range("D3").Value = "22:55:28"
$ws.Range("E3").Value = "22:55:28"
$ws.Range("F3").Value = "10:50:49"
$ws.Range("G3").Value = "10:57:15"

# Row 4: Time Total
$ws.Range("B4").Value = "0:0:30"
$ws.Range("C4").Value = "0:4:38"
$ws.Range("D4").Value = "13:0:0"
$ws.Range("E4").Value = "14:11:6"
$ws.Range("F4").Value = "0:5:39"
$ws.Range("G4").Value = "0:1:53"

# Row 5: Total Vehicles
$ws.Range("B5").Value = 19
$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 103
$ws.Range("G5").Value = 42

# Row 6: Vehicles from Left
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 39
$ws.Range("G6").Value = 20

# Row 7: Vehicles from Right
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 55
$ws.Range("D7").Value = 13
$ws.Range("E7").Value = 17
$ws.Range("F7").Value = 64
$ws.Range("G7").Value = 22

# Row 8: Vehicles from N/A
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# Row 9: Total Hours
$ws.Range("B9").Value = 0.008330000000000001
$ws.Range("C9").Value = 0.07722
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 14.185
$ws.Range("F9").Value = 0.09417
$ws.Range("G9").Value = 0.03139

# Row 10: Vehicles per Hour
$ws.Range("B10").Value = 2280.91
$ws.Range("C10").Value = 906.5
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1.27
$ws.Range("F10").Value = 1093.77
$ws.Range("G10").Value = 1338.01

# Row 11: Vehicles per Hour from Left
$ws.Range("B11").Value = 1320.53
$ws.Range("C11").Value = 194.25
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.07000000000000001
$ws.Range("F11").Value = 414.14
$ws.Range("G11").Value = 637.15

# Row 12: Vehicles per Hour from Right
$ws.Range("B12").Value = 960.38
$ws.Range("C12").Value = 712.25
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1.2
$ws.Range("F12").Value = 679.62
$ws.Range("G12").Value = 700.86

# Row 13: Morning Peak Start Time
$ws.Range("B13").Value = "N/A"
$ws.Range("C13").Value = "N/A"
$ws.Range("D13").Value = "N/A"
$ws.Range("E13").Value = "8:44:22"
$ws.Range("F13").Value = "N/A"
$ws.Range("G13").Value = "N/A"

# Row 14: Morning Peak End Time
$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "10:25:28"
$ws.Range("F14").Value = "N/A"
$ws.Range("G14").Value = "N/A"

# Row 15: Morning Peak Total Vehicles
$ws.Range("B15").Value = "N/A"
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = "N/A"
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = "N/A"
$ws.Range("G15").Value = "N/A"

# Row 16: Morning Peak Total Hours
$ws.Range("B16").Value = "N/A"
$ws.Range("C16").Value = "N/A"
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = 1.685
$ws.Range("F16").Value = "N/A"
$ws.Range("G16").Value = "N/A"

# Row 17: Morning Peak Vehicles per Hour
$ws.Range("B17").Value = "N/A"
$ws.Range("C17").Value = "N/A"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = 4.15
$ws.Range("F17").Value = "N/A"
$ws.Range("G17").Value = "N/A"

# Row 18: Morning Peak Vehicles from Left
$ws.Range("B18").Value = "N/A"
$ws.Range("C18").Value = "N/A"
$ws.Range("D18").Value = "N/A"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = "N/A"
$ws.Range("G18").Value = "N/A"

# Row 19: Morning Peak Vehicles from Right
$ws.Range("B19").Value = "N/A"
$ws.Range("C19").Value = "N/A"
$ws.Range("D19").Value = "N/A"
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = "N/A"
$ws.Range("G19").Value = "N/A"

# Row 20: Morning Peak VpH from Left
$ws.Range("B20").Value = "N/A"
$ws.Range("C20").Value = "N/A"
$ws.Range("D20").Value = "N/A"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = "N/A"
$ws.Range("G20").Value = "N/A"

# Row 21: Morning Peak VpH from Right
$ws.Range("B21").Value = "N/A"
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = "N/A"
$ws.Range("E21").Value = 4.15
$ws.Range("F21").Value = "N/A"
$ws.Range("G21").Value = "N/A"

# Row 22: Night Peak Start Time
$ws.Range("B22").Value = "N/A"
$ws.Range("C22").Value = "N/A"
$ws.Range("D22").Value = "17:55:28"
$ws.Range("E22").Value = "17:55:28"
$ws.Range("F22").Value = "N/A"
$ws.Range("G22").Value = "N/A"

# Row 23: Night Peak End Time
$ws.Range("B23").Value = "N/A"
$ws.Range("C23").Value = "N/A"
$ws.Range("D23").Value = "21:12:28"
$ws.Range("E23").Value = "21:12:28"
$ws.Range("F23").Value = "N/A"
$ws.Range("G23").Value = "N/A"

# Row 24: Night Peak Total Vehicles
$ws.Range("B24").Value = "N/A"
$ws.Range("C24").Value = "N/A"
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 6
$ws.Range("F24").Value = "N/A"
$ws.Range("G24").Value = "N/A"

# Row 25: Night Peak Total Hours
$ws.Range("B25").Value = "N/A"
$ws.Range("C25").Value = "N/A"
$ws.Range("D25").Value = 3.28333
$ws.Range("E25").Value = 3.28333
$ws.Range("F25").Value = "N/A"
$ws.Range("G25").Value = "N/A"

# Row 26: Night Peak Vehicles per Hour
$ws.Range("B26").Value = "N/A"
$ws.Range("C26").Value = "N/A"
$ws.Range("D26").Value = 1.83
$ws.Range("E26").Value = 1.83
$ws.Range("F26").Value = "N/A"
$ws.Range("G26").Value = "N/A"

# Row 27: Night Peak Vehicles from Left
$ws.Range("B27").Value = "N/A"
$ws.Range("C27").Value = "N/A"
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = "N/A"
$ws.Range("G27").Value = "N/A"

# Row 28: Night Peak Vehicles from Right
$ws.Range("B28").Value = "N/A"
$ws.Range("C28").Value = "N/A"
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 6
$ws.Range("F28").Value = "N/A"
$ws.Range("G28").Value = "N/A"

# Row 29: Night Peak Vph from Left
$ws.Range("B29").Value = "N/A"
$ws.Range("C29").Value = "N/A"
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = "N/A"
$ws.Range("G29").Value = "N/A"

# Row 30: Night Peak Vph from Right
$ws.Range("B30").Value = "N/A"
$ws.Range("C30").Value = "N/A"
$ws.Range("D30").Value = 1.83
$ws.Range("E30").Value = 1.83
$ws.Range("F30").Value = "N/A"
$ws.Range("G30").Value = "N/A"
